$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cells whose new value looks like a plain number but must stay text
# (matches the original inlineStr cell type in the source data).
$textCells = @('D5','D6','D9','D13','D17','D19','D20','D21','D22','D23','D24','D25','D26','D28','D35','D36','D37','D41','D42','D45','D46','D48')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = '@'
}

$ws.Range('D2').Value = '67.228.27'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '3.111.58'
$ws.Range('E3').Value = '  -0.76%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '579.31'
$ws.Range('E5').Value = '  -0.23%  '
$ws.Range('D6').Value = '172.87'
$ws.Range('E6').Value = '  -1.05%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.99%  '
$ws.Range('D9').Value = '6.51'
$ws.Range('E9').Value = '  +0.26%  '
$ws.Range('E10').Value = '  -1.66%  '
$ws.Range('E11').Value = '  -1.74%  '
$ws.Range('D13').Value = '36.68'
$ws.Range('E13').Value = '  -1.72%  '
$ws.Range('E14').Value = '  -2.06%  '
$ws.Range('D15').Value = '3.626.92'
$ws.Range('E15').Value = '  -0.57%  '
$ws.Range('D16').Value = '67.128.74'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').Value = '7.09'
$ws.Range('E17').Value = '  -1.51%  '
$ws.Range('D18').Value = '3.109.11'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').Value = '16.63'
$ws.Range('E19').Value = '  +2.33%  '
$ws.Range('D20').Value = '490.34'
$ws.Range('E20').Value = '  +0.41%  '
$ws.Range('B21').Value = 'Polygon'
$ws.Range('C21').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D21').Value = '0.699'
$ws.Range('E21').Value = '  -2.74%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '7.81'
$ws.Range('E22').Value = '  +1.15%  '
$ws.Range('D23').Value = '83.82'
$ws.Range('E23').Value = '  -0.75%  '
$ws.Range('D24').Value = '13.09'
$ws.Range('E24').Value = '  -2.22%  '
$ws.Range('D25').Value = '2.29'
$ws.Range('E25').Value = '  -2.22%  '
$ws.Range('D26').Value = '10.55'
$ws.Range('E26').Value = '  +4.45%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').Value = '7.88'
$ws.Range('E28').Value = '  -1.96%  '
$ws.Range('E29').Value = '  -2.68%  '
$ws.Range('E30').Value = '  -1.13%  '
$ws.Range('E31').Value = '  -2.71%  '
$ws.Range('E32').Value = '  -1.33%  '
$ws.Range('D33').Value = '0.0₃0941'
$ws.Range('E33').Value = '  -6.57%  '
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('D35').Value = '5.78'
$ws.Range('E35').Value = '  -3.19%  '
$ws.Range('D36').Value = '0.970'
$ws.Range('E36').Value = '  -2.35%  '
$ws.Range('D37').Value = '46.76'
$ws.Range('E37').Value = '  -1.35%  '
$ws.Range('E38').Value = '  -4.53%  '
$ws.Range('E39').Value = '  +0.23%  '
$ws.Range('E40').Value = '  -2.47%  '
$ws.Range('D41').Value = '8.45'
$ws.Range('E41').Value = '  -2.70%  '
$ws.Range('D42').Value = '386.13'
$ws.Range('E42').Value = '  -0.73%  '
$ws.Range('D43').Value = '2.800.45'
$ws.Range('E43').Value = '  -2.20%  '
$ws.Range('E44').Value = '  -8.96%  '
$ws.Range('D45').Value = '0.0348'
$ws.Range('E45').Value = '  -3.21%  '
$ws.Range('D46').Value = '135.25'
$ws.Range('E46').Value = '  -0.83%  '
$ws.Range('D48').Value = '24.97'
$ws.Range('E48').Value = '  -0.94%  '
$ws.Range('E49').Value = '  -2.04%  '
$ws.Range('E50').Value = '  -1.94%  '
$ws.Range('E51').Value = '  -1.81%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = 'Normal'
}
